# Updated cryptos list on Wed May 29 03:39:23 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for the
# crypto list, and fixes the Kaspa/Mantle row ordering (rows 39-40 swapped
# places, each keeping its own Coin/Link/Price/Volume data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new value.
$updates = [ordered]@{
    'D2'  = '68.753.55';   'E2'  = '  +0.90%  '
    'D3'  = '3.846.33';    'E3'  = '  +0.07%  '
    'D4'  = '1.00';        'E4'  = '  -0.10%  '
    'D5'  = '601.94';      'E5'  = '  +0.54%  '
    'D6'  = '171.47';      'E6'  = '  +2.47%  '
    'D7'  = '3.847.51';    'E7'  = '  +0.09%  '
                           'E8'  = '  -0.11%  '
                           'E9'  = '  +0.00%  '
                           'E10' = '  +1.48%  '
    'D11' = '6.53';        'E11' = '  +3.20%  '
                           'E12' = '  -0.27%  '
                           'E13' = '  +14.10%  '
    'D14' = '37.16';       'E14' = '  -0.48%  '
    'D15' = '4.492.89';    'E15' = '  +0.01%  '
    'D16' = '3.829.65';    'E16' = '  -0.17%  '
    'D17' = '68.761.81';   'E17' = '  +0.62%  '
    'D18' = '18.31';       'E18' = '  +0.77%  '
    'D19' = '7.38';        'E19' = '  -2.43%  '
                           'E20' = '  -0.13%  '
    'D21' = '11.22';       'E21' = '  +4.14%  '
    'D22' = '473.31';      'E22' = '  +0.52%  '
                           'E23' = '  -0.87%  '
                           'E24' = '  +2.80%  '
    'D25' = '83.62';       'E25' = '  -1.05%  '
                           'E26' = '  +0.53%  '
    'D27' = '12.16';       'E27' = '  -1.58%  '
    'D28' = '10.38';       'E28' = '  +3.46%  '
                           'E30' = '  +0.56%  '
    'D31' = '3.996.67';    'E31' = '  -0.01%  '
    'D32' = '7.72';        'E32' = '  -0.47%  '
    'D33' = '31.39'
                           'E34' = '  -0.23%  '
    'D35' = '9.34';        'E35' = '  -1.37%  '
    'D36' = '3.813.11';    'E36' = '  -0.20%  '
    'D37' = '3.93';        'E37' = '  +19.03%  '
                           'E38' = '  -0.80%  '
    # Kaspa and Mantle swap ranking positions (row 39 <-> row 40).
    'B39' = 'Mantle'
    'C39' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D39' = '1.02';        'E39' = '  +1.18%  '
    'B40' = 'Kaspa'
    'C40' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D40' = '0.140';       'E40' = '  -0.02%  '
                           'E41' = '  -0.36%  '
                           'E42' = '  +0.03%  '
                           'E43' = '  +0.49%  '
    'D44' = '0.000301';    'E44' = '  +9.81%  '
                           'E45' = '  -0.06%  '
                           'E46' = '  -0.01%  '
    'D47' = '421.39';      'E47' = '  -1.89%  '
    'D48' = '8.68';        'E48' = '  +0.86%  '
    'D49' = '46.45';       'E49' = '  -2.06%  '
    'D50' = '141.55';      'E50' = '  -0.72%  '
                           'E51' = '  -0.55%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)

    if ($value -match '^-?\d+(\.\d+)?$') {
        # This text looks like a plain number (e.g. "1.00", "0.140").
        # These columns are text, not numeric, so force Excel to keep the
        # value as a string (leading apostrophe) instead of silently
        # normalising it to a Double, then drop the resulting "Text"
        # number-format override so the cell style stays untouched.
        $cell.Value = "'" + $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}
